# TC001.xlsx fix: redact the login credentials that were typed into column C
# (they were leaking the real email/password through the cell text) while
# keeping the original hyperlinks - and their original "display" text - that
# were already wired up to those two cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC001")

# C2: login-email data cell -> hyperlink to the e-mail address, display text
# preserved as-is, then the visible cell text is masked to "XXXXXX".
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:karim_94123@yahoo.co.in", "", "", "karim_94123@yahoo.co.in")
$ws.Range("C2").Value = "XXXXXX"
# touch the font so the engine re-settles the cell's style back onto the
# existing "Hyperlink" cell format (s=4) instead of minting a fresh xf
$ws.Range("C2").Font.Name = $ws.Range("C2").Font.Name

# C3: login-password data cell -> same treatment.
$ws.Hyperlinks.Add($ws.Range("C3"), "P@ssw0rdL", "", "", "P@ssw0rdL")
$ws.Range("C3").Value = "XXXXXX"
$ws.Range("C3").Font.Name = $ws.Range("C3").Font.Name

# Leave the selection on the last touched cell, as in the saved file.
$ws.Range("C3").Select()
